$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates: row -> new value.
# These cells are plain-text price strings (e.g. "49.567.25", "112.48").
# Excel auto-detects numeric-looking text as a Number, so we force the
# cell to Text format before assigning, then strip the formatting change
# back off (ClearFormats) so the cell keeps its original (default) style
# while the stored value remains a text string.
$priceUpdates = @{
    2  = "49.567.25"
    3  = "2.645.57"
    5  = "112.48"
    6  = "326.61"
    8  = "1.00"
    10 = "39.66"
    11 = "19.97"
    14 = "7.63"
    15 = "3.061.43"
    16 = "2.639.44"
    18 = "49.573.26"
    19 = "13.39"
    23 = "268.54"
    24 = "68.95"
    27 = "1.00"
    28 = "10.16"
    31 = "34.63"
    33 = "5.48"
    35 = "19.18"
    38 = "2.03"
    40 = "129.97"
    41 = "23.52"
    42 = "0.0348"
    43 = "2.29"
    45 = "2.062.55"
    47 = "2.10"
    50 = "5.25"
    51 = "58.86"
}

# Column E (Volume 1h) updates: row -> new value (keeping the "  +x.xx%  " / "  -x.xx%  " padding)
$volumeUpdates = @{
    2  = "  -0.81%  "
    3  = "  -0.06%  "
    4  = "  +0.07%  "
    5  = "  -1.25%  "
    6  = "  -0.07%  "
    7  = "  -1.09%  "
    8  = "  +0.00%  "
    9  = "  -1.50%  "
    10 = "  -3.41%  "
    11 = "  -1.08%  "
    12 = "  -0.88%  "
    13 = "  +1.87%  "
    14 = "  +3.16%  "
    15 = "  +0.02%  "
    16 = "  -0.47%  "
    17 = "  -1.67%  "
    18 = "  -0.62%  "
    19 = "  +1.42%  "
    20 = "  -1.65%  "
    21 = "  -0.43%  "
    22 = "  -0.96%  "
    23 = "  -2.95%  "
    24 = "  -4.35%  "
    25 = "  -0.90%  "
    26 = "  -2.77%  "
    27 = "  +0.05%  "
    28 = "  +1.51%  "
    29 = "  -1.00%  "
    30 = "  -2.43%  "
    31 = "  -4.26%  "
    32 = "  -1.33%  "
    33 = "  +0.31%  "
    34 = "  +1.52%  "
    35 = "  -1.48%  "
    36 = "  -0.12%  "
    37 = "  -1.89%  "
    38 = "  -1.65%  "
    39 = "  +0.36%  "
    40 = "  +4.79%  "
    41 = "  +6.28%  "
    42 = "  +10.62%  "
    43 = "  +2.87%  "
    45 = "  -1.05%  "
    46 = "  -0.76%  "
    47 = "  +5.78%  "
    48 = "  -4.79%  "
    49 = "  -2.31%  "
    50 = "  -2.89%  "
    51 = "  -1.50%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.ClearFormats()
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
